# Auto-generated edit script: applies the cryptos.xlsx price/volume refresh
# described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.874.91'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '3.207.93'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.02'
$ws.Range("E5").Value = '  +4.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.25'
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.205.68'
$ws.Range("E8").Value = '  +0.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  +1.40%  '
$ws.Range("E10").Value = '  -1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.13'
$ws.Range("E11").Value = '  -1.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.509'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").Value = '  -2.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.43'
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("D15").Value = '3.734.13'
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").Value = '65.989.55'
$ws.Range("E16").Value = '  +1.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.37'
$ws.Range("E17").Value = '  +2.28%  '
$ws.Range("D18").Value = '3.208.75'
$ws.Range("E18").Value = '  +0.72%  '
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '510.70'
$ws.Range("E20").Value = '  -0.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.78'
$ws.Range("E21").Value = '  +5.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.735'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.12'
$ws.Range("E23").Value = '  -2.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.96'
$ws.Range("E24").Value = '  +1.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.27'
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.25'
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("E28").Value = '  +3.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +1.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.86'
$ws.Range("E30").Value = '  +2.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.79'
$ws.Range("E31").Value = '  +7.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.08'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.63'
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.43'
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("E37").Value = '  +2.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '486.53'
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0421'
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.99'
$ws.Range("E40").Value = '  -4.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.85'
$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("D42").Value = '3.024.51'
$ws.Range("E42").Value = '  -3.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.119'
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.292'
$ws.Range("E44").Value = '  +1.33%  '
$ws.Range("D45").Value = '0.0₃0643'
$ws.Range("E45").Value = '  +6.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.45'
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '29.04'
$ws.Range("E47").Value = '  -1.56%  '
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.32'
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.02'
$ws.Range("E51").Value = '  -1.35%  '
